# Cyclic shift of the observation rows 2-8: each row's data (for the
# columns that actually differ between rows) moves up by one row, and the
# original row 2 data wraps around to row 8. Row numbers (column A's
# underlying "Id") travel with everything else because they are part of
# the same record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "before" values for every column that varies row-to-row,
# keyed by (column, row).
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "Y", "AA", "AX")

$before = @{}
foreach ($col in $cols) {
    for ($r = 2; $r -le 8; $r++) {
        $before[$col + $r] = $ws.Range($col + $r).Value2
    }
}

function Set-CellValue($range, $value) {
    # These source cells are plain text (e.g. ISO-looking dates stored as
    # literal strings). Writing such a string back through COM would make
    # Excel auto-convert it to a real date serial, which the original file
    # does not use here. A leading apostrophe forces text entry, exactly
    # like typing it in the UI, and is not retained in the stored value.
    if ($value -is [string] -and $value -match '^\d{4}-\d{2}-\d{2}$') {
        $range.Value2 = "'" + $value
    } else {
        $range.Value2 = $value
    }
}

# new row r (2..7) <- old row r+1 ; new row 8 <- old row 2
for ($r = 2; $r -le 7; $r++) {
    $srcRow = $r + 1
    foreach ($col in $cols) {
        Set-CellValue $ws.Range($col + $r) $before[$col + $srcRow]
    }
}
foreach ($col in $cols) {
    Set-CellValue $ws.Range($col + "8") $before[$col + "2"]
}
